# Rewrite the lookup-list ranges on "Dmstc Stndrd Upld Tmplt" (sheet2):
# columns B..M and AC hold sequential "prefix:N" helper values used for
# data-validation / lookup lists. Each column's list is being
# lengthened or shortened, so we clear the old ranges first and then
# write the new, correctly-sized lists back in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dmstc Stndrd Upld Tmplt")

# Clear out the previous lookup-list content (rows 6-52 for B:M, rows
# 6-52 for AC) so stale rows beyond the new list lengths disappear
# completely instead of lingering as blank styled cells.
$ws.Range("B6:M52").Clear()
$ws.Range("AC6:AC52").Clear()

# column letter -> (prefix, item count) for the new lists
$lists = @{
    "B"  = @("vpn", 17)
    "C"  = @("ven", 6)
    "D"  = @("lab", 30)
    "E"  = @("col", 33)
    "F"  = @("siz", 11)
    "G"  = @("typ", 44)
    "H"  = @("pon", 23)
    "I"  = @("sup", 3)
    "J"  = @("dep", 4)
    "K"  = @("cla", 13)
    "L"  = @("cat", 31)
    "M"  = @("mat", 41)
}

foreach ($col in $lists.Keys) {
    $prefix = $lists[$col][0]
    $count = $lists[$col][1]
    for ($i = 0; $i -lt $count; $i++) {
        $row = 6 + $i
        $ws.Range($col + $row).Value = $prefix + ":" + $i
    }
}

# AC's list starts one row earlier (row 5, inside the hidden header
# row) and has 42 items (wei:0 .. wei:41).
$wei_count = 42
for ($i = 0; $i -lt $wei_count; $i++) {
    $row = 5 + $i
    $ws.Range("AC" + $row).Value = "wei:" + $i
}
